$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '48.145.28'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.17%  '
$ws.Range('E2').ClearFormats()

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.518.74'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.24%  '
$ws.Range('E3').ClearFormats()

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('E4').ClearFormats()

$ws.Range('B5').NumberFormat = '@'
$ws.Range('B5').Value = 'BNB'
$ws.Range('B5').ClearFormats()
$ws.Range('C5').NumberFormat = '@'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('C5').ClearFormats()
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.13'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.52%  '
$ws.Range('E5').ClearFormats()

$ws.Range('B6').NumberFormat = '@'
$ws.Range('B6').Value = 'Solana'
$ws.Range('B6').ClearFormats()
$ws.Range('C6').NumberFormat = '@'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('C6').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '109.35'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.94%  '
$ws.Range('E6').ClearFormats()

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +2.06%  '
$ws.Range('E7').ClearFormats()

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.98%  '
$ws.Range('E9').ClearFormats()

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.47'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.53%  '
$ws.Range('E10').ClearFormats()

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.54'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +12.88%  '
$ws.Range('E11').ClearFormats()

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.06%  '
$ws.Range('E12').ClearFormats()

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.15%  '
$ws.Range('E13').ClearFormats()

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.27'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.09%  '
$ws.Range('E14').ClearFormats()

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.913.82'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.21%  '
$ws.Range('E15').ClearFormats()

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.517.26'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.12%  '
$ws.Range('E16').ClearFormats()

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.08%  '
$ws.Range('E17').ClearFormats()

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '47.981.43'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.00%  '
$ws.Range('E18').ClearFormats()

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.23'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +4.33%  '
$ws.Range('E19').ClearFormats()

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.64%  '
$ws.Range('E20').ClearFormats()

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.54%  '
$ws.Range('E21').ClearFormats()

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.70'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.08%  '
$ws.Range('E22').ClearFormats()

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.09'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.01%  '
$ws.Range('E23').ClearFormats()

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '264.10'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +7.60%  '
$ws.Range('E24').ClearFormats()

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.57'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.82%  '
$ws.Range('E25').ClearFormats()

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('E26').ClearFormats()

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '26.09'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.68%  '
$ws.Range('E27').ClearFormats()

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.13'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.93%  '
$ws.Range('E28').ClearFormats()

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.46%  '
$ws.Range('E29').ClearFormats()

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.13%  '
$ws.Range('E30').ClearFormats()

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '36.31'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +3.74%  '
$ws.Range('E31').ClearFormats()

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.38%  '
$ws.Range('E32').ClearFormats()

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.73'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.51%  '
$ws.Range('E33').ClearFormats()

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.07%  '
$ws.Range('E34').ClearFormats()

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E35').ClearFormats()

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0790'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.12%  '
$ws.Range('E36').ClearFormats()

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.99'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.58%  '
$ws.Range('E37').ClearFormats()

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.41%  '
$ws.Range('E38').ClearFormats()

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.99'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.66%  '
$ws.Range('E39').ClearFormats()

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('E40').ClearFormats()

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.85%  '
$ws.Range('E41').ClearFormats()

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '119.84'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.81%  '
$ws.Range('E42').ClearFormats()

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.87%  '
$ws.Range('E43').ClearFormats()

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0301'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.49%  '
$ws.Range('E44').ClearFormats()

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.017.06'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.93%  '
$ws.Range('E45').ClearFormats()

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +5.30%  '
$ws.Range('E46').ClearFormats()

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.93'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +9.46%  '
$ws.Range('E47').ClearFormats()

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.88%  '
$ws.Range('E48').ClearFormats()

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.87%  '
$ws.Range('E49').ClearFormats()

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.23'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.07%  '
$ws.Range('E50').ClearFormats()

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.91'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.75%  '
$ws.Range('E51').ClearFormats()
